$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Copy header style from D1 (existing header) to E1:F1
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values
$ws.Range("E2").Value = 4.036300000734627
$ws.Range("F2").Value = 0

$ws.Range("E3").Value = 3.871199995046481
$ws.Range("F3").Value = 4096

$ws.Range("E4").Value = 27.28549999301322
$ws.Range("F4").Value = 0

$ws.Range("E5").Value = 1.443799992557615
$ws.Range("F5").Value = 0
